$wb = $excel.ActiveWorkbook

# ALC row 4 (item id 5470)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 280
$ws.Range("I4").Value = 280
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 280
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -166
$ws.Range("N4").ClearContents()

# ALC row 33 (item id 5512)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 255.22858
$ws.Range("I33").Value = 237.44827
$ws.Range("K33").Value = 237.44827
$ws.Range("M33").Value = -8.448270000000008

# ALC row 52 (item id 4567)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H52").Value = 298797
$ws.Range("I52").Value = 200198
$ws.Range("J52").Value = 594594
$ws.Range("K52").Value = 600594
$ws.Range("L52").Value = 1783782
$ws.Range("M52").Value = -600434
$ws.Range("N52").Value = -1784102

# ALC row 53 (item id 5479)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 454.64285
$ws.Range("I53").Value = 265.7143
$ws.Range("J53").Value = 643.5714
$ws.Range("K53").Value = 265.7143
$ws.Range("L53").Value = 643.5714
$ws.Range("M53").Value = 371.2857
$ws.Range("N53").Value = -1917.5714

# ALC row 62 (item id 27781)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3315.8572
$ws.Range("I62").Value = 3801
$ws.Range("J62").Value = 2103
$ws.Range("K62").Value = 3801
$ws.Range("L62").Value = 2103
$ws.Range("M62").Value = -3177
$ws.Range("N62").Value = -3351

# ALC row 65 (item id 27781)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 3315.8572
$ws.Range("I65").Value = 3801
$ws.Range("J65").Value = 2103
$ws.Range("K65").Value = 19005
$ws.Range("L65").Value = 10515
$ws.Range("M65").Value = -15885
$ws.Range("N65").Value = -16755

# ALC row 116 (item id 27778)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 1898.3334
$ws.Range("I116").Value = 1798
$ws.Range("J116").Value = 2400
$ws.Range("K116").Value = 1798
$ws.Range("L116").Value = 2400
$ws.Range("M116").Value = 1644
$ws.Range("N116").Value = -9284

# ALC row 129 (item id 36115)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 829.7941
$ws.Range("J129").Value = 1019.75
$ws.Range("L129").Value = 3059.25
$ws.Range("N129").Value = -13059.25

# ALC row 132 (item id 44049)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 756943.1
$ws.Range("I132").Value = 2877.3396
$ws.Range("K132").Value = 8632.0188
$ws.Range("M132").Value = -6102.0188

# ALC row 139 (item id 42306)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H139").Value = 38000
$ws.Range("J139").Value = 38000
$ws.Range("L139").Value = 38000
$ws.Range("N139").Value = -48280

# ARM row 132 (item id 43997)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 9846824
$ws.Range("I132").Value = 11652792
$ws.Range("J132").Value = 139747.25
$ws.Range("K132").Value = 34958376
$ws.Range("L132").Value = 419241.75
$ws.Range("M132").Value = -34955846
$ws.Range("N132").Value = -424301.75

# BSM row 134 (item id 43998)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 10755746
$ws.Range("I134").Value = 2968.4707
$ws.Range("J134").Value = 23812692
$ws.Range("K134").Value = 8905.4121
$ws.Range("L134").Value = 71438076
$ws.Range("M134").Value = -6370.4121
$ws.Range("N134").Value = -71443146

# CRP row 94 (item id 32934)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 3386.04
$ws.Range("I94").Value = 7194.6665
$ws.Range("K94").Value = 7194.6665
$ws.Range("M94").Value = -6743.6665

# CUL row 34 (item id 4749)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 633.6667
$ws.Range("I34").Value = 300
$ws.Range("J34").Value = 900.6
$ws.Range("K34").Value = 900
$ws.Range("L34").Value = 2701.8
$ws.Range("M34").Value = -816
$ws.Range("N34").Value = -2869.8

# CUL row 63 (item id 12866)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 4189.857
$ws.Range("I63").Value = 1733.3334
$ws.Range("J63").Value = 4859.8184
$ws.Range("K63").Value = 5200.0002
$ws.Range("L63").Value = 14579.4552
$ws.Range("M63").Value = -4451.0002
$ws.Range("N63").Value = -16077.4552

# CUL row 66 (item id 12866)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H66").Value = 4189.857
$ws.Range("I66").Value = 1733.3334
$ws.Range("J66").Value = 4859.8184
$ws.Range("K66").Value = 15600.0006
$ws.Range("L66").Value = 43738.3656
$ws.Range("M66").Value = -11856.0006
$ws.Range("N66").Value = -51226.3656

# CUL row 104 (item id 19807)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H104").Value = 3997.7778
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 3997.7778
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 11993.3334
$ws.Range("M104").ClearContents()
$ws.Range("N104").Value = -17235.3334

# CUL row 129 (item id 36054)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 3474436.5
$ws.Range("I129").Value = 1858.5714
$ws.Range("J129").Value = 4904321.5
$ws.Range("K129").Value = 5575.7142
$ws.Range("L129").Value = 14712964.5
$ws.Range("M129").Value = -575.7142000000003
$ws.Range("N129").Value = -14722964.5

# CUL row 131 (item id 36060)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 864.53125
$ws.Range("I131").Value = 480.9
$ws.Range("J131").Value = 935.5741
$ws.Range("K131").Value = 1442.7
$ws.Range("L131").Value = 2806.7223
$ws.Range("M131").Value = 3597.3
$ws.Range("N131").Value = -12886.7223

# GSM row 102 (item id 36169)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3518
$ws.Range("I102").Value = 4937.3335
$ws.Range("J102").Value = 2453.5
$ws.Range("K102").Value = 4937.3335
$ws.Range("L102").Value = 2453.5
$ws.Range("M102").Value = -3315.3335
$ws.Range("N102").Value = -5697.5

# GSM row 122 (item id 36182)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1877.1464
$ws.Range("I122").Value = 1575.1562
$ws.Range("J122").Value = 2950.889
$ws.Range("K122").Value = 4725.4686
$ws.Range("L122").Value = 8852.667000000001
$ws.Range("M122").Value = -2275.4686
$ws.Range("N122").Value = -13752.667

# GSM row 123 (item id 34150)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 23090.75
$ws.Range("J123").Value = 23090.75
$ws.Range("L123").Value = 23090.75
$ws.Range("N123").Value = -27990.75

# GSM row 124 (item id 34247)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H124").Value = 58900
$ws.Range("J124").Value = 58900
$ws.Range("L124").Value = 58900
$ws.Range("N124").Value = -68720

# LTW row 22 (item id 5277)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1065.0588
$ws.Range("I22").Value = 1018.75
$ws.Range("J22").Value = 1106.2222
$ws.Range("K22").Value = 1018.75
$ws.Range("L22").Value = 1106.2222
$ws.Range("M22").Value = -723.75
$ws.Range("N22").Value = -1696.2222

# LTW row 27 (item id 5277)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 1065.0588
$ws.Range("I27").Value = 1018.75
$ws.Range("J27").Value = 1106.2222
$ws.Range("K27").Value = 1018.75
$ws.Range("L27").Value = 1106.2222
$ws.Range("M27").Value = -911.75
$ws.Range("N27").Value = -1320.2222

# LTW row 40 (item id 36248)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2831.3076
$ws.Range("I40").Value = 2686
$ws.Range("K40").Value = 2686
$ws.Range("M40").Value = -2550

# LTW row 61 (item id 27740)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1949.409
$ws.Range("I61").Value = 1918.8
$ws.Range("J61").Value = 2015
$ws.Range("K61").Value = 1918.8
$ws.Range("L61").Value = 2015
$ws.Range("M61").Value = -1716.8
$ws.Range("N61").Value = -2419

# LTW row 82 (item id 12565)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2243.3333
$ws.Range("I82").Value = 1275.5
$ws.Range("J82").Value = 2727.25
$ws.Range("K82").Value = 1275.5
$ws.Range("L82").Value = 2727.25
$ws.Range("M82").Value = -914.5
$ws.Range("N82").Value = -3449.25

# LTW row 85 (item id 12565)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 2243.3333
$ws.Range("I85").Value = 1275.5
$ws.Range("J85").Value = 2727.25
$ws.Range("K85").Value = 1275.5
$ws.Range("L85").Value = 2727.25
$ws.Range("M85").Value = -27.5
$ws.Range("N85").Value = -5223.25

# LTW row 93 (item id 19993)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1442.125
$ws.Range("I93").Value = 1497.5
$ws.Range("J93").Value = 1276
$ws.Range("K93").Value = 1497.5
$ws.Range("L93").Value = 1276
$ws.Range("M93").Value = -249.5
$ws.Range("N93").Value = -3772

# LTW row 100 (item id 19995)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 1893.9375
$ws.Range("I100").Value = 1580.6
$ws.Range("K100").Value = 1580.6
$ws.Range("M100").Value = -1039.6

# LTW row 113 (item id 27740)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 1949.409
$ws.Range("I113").Value = 1918.8
$ws.Range("J113").Value = 2015
$ws.Range("K113").Value = 1918.8
$ws.Range("L113").Value = 2015
$ws.Range("M113").Value = 251.2
$ws.Range("N113").Value = -6355

# LTW row 122 (item id 36247)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3942.0667
$ws.Range("I122").Value = 4092
$ws.Range("J122").Value = 3770.7144
$ws.Range("K122").Value = 12276
$ws.Range("L122").Value = 11312.1432
$ws.Range("M122").Value = -9826
$ws.Range("N122").Value = -16212.1432

# LTW row 123 (item id 35408)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H123").Value = 250018750
$ws.Range("J123").Value = 250018750
$ws.Range("L123").Value = 250018750
$ws.Range("N123").Value = -250028550

# LTW row 132 (item id 44058)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 39472.93
$ws.Range("I132").Value = 4538.1875
$ws.Range("J132").Value = 82469.53999999999
$ws.Range("K132").Value = 13614.5625
$ws.Range("L132").Value = 247408.62
$ws.Range("M132").Value = -11084.5625
$ws.Range("N132").Value = -252468.62

# WVR row 107 (item id 27746)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 534.9259
$ws.Range("I107").Value = 407.33334
$ws.Range("J107").Value = 981.5
$ws.Range("K107").Value = 1222.00002
$ws.Range("L107").Value = 2944.5
$ws.Range("M107").Value = 697.9999800000001
$ws.Range("N107").Value = -6784.5

# WVR row 122 (item id 36208)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1808.9166
$ws.Range("I122").Value = 1501.4615
$ws.Range("K122").Value = 4504.3845
$ws.Range("M122").Value = -2054.3845
